$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (pushes old rows 54..81 down to 55..82)
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the same static columns as its neighbours,
# and the new data values from the diff.
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 44755
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = 100112012
$ws.Cells.Item(54, 7).Value = "Espinaca"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 40
$ws.Cells.Item(54, 11).Value = 9000
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 9500
$ws.Cells.Item(54, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 950
$ws.Cells.Item(54, 17).Value = 10
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows in column D
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
